$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.727.00'
$ws.Range("E2").Value = '  +2.93%  '
$ws.Range("D3").Value = '1.862.06'
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.039'
$ws.Range("E4").Value = '  +3.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.74'
$ws.Range("E5").Value = '  +3.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.034'
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4413'
$ws.Range("E7").Value = '  +2.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3803'
$ws.Range("E8").Value = '  +2.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07448'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8822'
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.72'
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("D12").Value = '1.872.25'
$ws.Range("E12").Value = '  -8.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.555'
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.738'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07205'
$ws.Range("E15").Value = '  +4.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.64'
$ws.Range("E16").Value = '  +3.56%  '
$ws.Range("E17").Value = '  +2.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009085'
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.52'
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("D21").Value = '27.752.67'
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.301'
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.43'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.69'
$ws.Range("E24").Value = '  +2.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.936'
$ws.Range("E25").Value = '  +2.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.84'
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.988'
$ws.Range("E27").Value = '  +4.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.311'
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.64'
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09091'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("E31").Value = '  +4.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7662'
$ws.Range("E32").Value = '  +3.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.574'
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.886'
$ws.Range("E34").Value = '  +2.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.034'
$ws.Range("E35").Value = '  +2.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.160'
$ws.Range("E36").Value = '  +3.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01982'
$ws.Range("E37").Value = '  +3.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05337'
$ws.Range("E38").Value = '  +2.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5191'
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.834'
$ws.Range("E40").Value = '  +3.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1690'
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.828'
$ws.Range("E42").Value = '  +5.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.667'
$ws.Range("E43").Value = '  +4.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.51'
$ws.Range("E44").Value = '  +2.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.60'
$ws.Range("E45").Value = '  +2.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.722'
$ws.Range("E46").Value = '  +4.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4682'
$ws.Range("E47").Value = '  +2.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06429'
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.854'
$ws.Range("E49").Value = '  +3.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.69'
$ws.Range("E50").Value = '  +4.77%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9343'
$ws.Range("E51").Value = '  +2.60%  '
